$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.313.11'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.96%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.249.77'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.75'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.638'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '64.36'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.03%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0950'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -9.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.93'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '26.31'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.585.05'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.88'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -5.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.03'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.819'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.252.26'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '43.260.21'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0966'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -5.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.88'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.90%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '246.87'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.00%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +16.24%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.44'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '173.82'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.18%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '21.53'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.47%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.14%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.125'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.90'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0677'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.90'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.61'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.38'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.71%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0249'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.91%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.79'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +5.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.48'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.49%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.09'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.60%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '96.51'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.15'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.90%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0937'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.96%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.422.88'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.60%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.26%  '
